# Auto-generated Excel COM-interop script
# Applies updated Leve profit calculations to multiple job sheets
# (ALC, ARM, BSM, CRP, CUL, LTW, WVR) per scheduled market-price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 10870840
$ws.Range("I135").Value = 342.96667
$ws.Range("J135").Value = 31253020
$ws.Range("K135").Value = 3086.70003
$ws.Range("L135").Value = 281277180
$ws.Range("M135").Value = -551.70003
$ws.Range("N135").Value = -281282250
$ws.Range("H137").Value = 33394.355
$ws.Range("I137").Value = 59720.117
$ws.Range("J137").Value = 1427.3572
$ws.Range("K137").Value = 179160.351
$ws.Range("L137").Value = 4282.071599999999
$ws.Range("M137").Value = -176610.351
$ws.Range("N137").Value = -9382.071599999999
$ws.Range("H141").Value = 43528.684
$ws.Range("I141").Value = 3077
$ws.Range("J141").Value = 71533.69500000001
$ws.Range("K141").Value = 9231
$ws.Range("L141").Value = 214601.085
$ws.Range("M141").Value = -4051
$ws.Range("N141").Value = -224961.085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5718.65
$ws.Range("I2").Value = 774.38464
$ws.Range("J2").Value = 14900.857
$ws.Range("K2").Value = 774.38464
$ws.Range("L2").Value = 14900.857
$ws.Range("M2").Value = -661.38464
$ws.Range("N2").Value = -15126.857
$ws.Range("H45").Value = 1367.4286
$ws.Range("J45").Value = 1670.5
$ws.Range("L45").Value = 1670.5
$ws.Range("N45").Value = -2424.5
$ws.Range("H102").Value = 2472.8572
$ws.Range("I102").Value = 2468.3333
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2468.3333
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -846.3332999999998
$ws.Range("N102").Value = -5744
$ws.Range("H105").Value = 35235
$ws.Range("J105").Value = 35235
$ws.Range("L105").Value = 35235
$ws.Range("N105").Value = -42223
$ws.Range("H116").Value = 5718.65
$ws.Range("I116").Value = 774.38464
$ws.Range("J116").Value = 14900.857
$ws.Range("K116").Value = 774.38464
$ws.Range("L116").Value = 14900.857
$ws.Range("M116").Value = 1519.61536
$ws.Range("N116").Value = -19488.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5718.65
$ws.Range("I3").Value = 774.38464
$ws.Range("J3").Value = 14900.857
$ws.Range("K3").Value = 774.38464
$ws.Range("L3").Value = 14900.857
$ws.Range("M3").Value = -660.38464
$ws.Range("N3").Value = -15128.857
$ws.Range("H88").Value = 16311
$ws.Range("I88").Value = 16311
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 16311
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -15905
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 16311
$ws.Range("I91").Value = 16311
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 16311
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -14907
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 67695.25
$ws.Range("J92").Value = 67695.25
$ws.Range("L92").Value = 67695.25
$ws.Range("N92").Value = -72687.25
$ws.Range("H100").Value = 22500
$ws.Range("J100").Value = 22500
$ws.Range("L100").Value = 22500
$ws.Range("N100").Value = -24664
$ws.Range("H108").Value = 22750
$ws.Range("I108").Value = 16500
$ws.Range("J108").Value = 29000
$ws.Range("K108").Value = 16500
$ws.Range("L108").Value = 29000
$ws.Range("M108").Value = -12660
$ws.Range("N108").Value = -36680
$ws.Range("H110").Value = 65940.39999999999
$ws.Range("J110").Value = 65940.39999999999
$ws.Range("L110").Value = 65940.39999999999
$ws.Range("N110").Value = -74120.39999999999
$ws.Range("H111").Value = 44150.5
$ws.Range("J111").Value = 44150.5
$ws.Range("L111").Value = 44150.5
$ws.Range("N111").Value = -52330.5
$ws.Range("H112").Value = 33999
$ws.Range("J112").Value = 33999
$ws.Range("L112").Value = 33999
$ws.Range("N112").Value = -36953
$ws.Range("H116").Value = 23333
$ws.Range("J116").Value = 23333
$ws.Range("L116").Value = 23333
$ws.Range("N116").Value = -32511
$ws.Range("H119").Value = 29472
$ws.Range("J119").Value = 29472
$ws.Range("L119").Value = 29472
$ws.Range("N119").Value = -39148
$ws.Range("H120").Value = 29996.666
$ws.Range("J120").Value = 29996.666
$ws.Range("L120").Value = 29996.666
$ws.Range("N120").Value = -39672.666
$ws.Range("H124").Value = 35923.75
$ws.Range("J124").Value = 35923.75
$ws.Range("L124").Value = 35923.75
$ws.Range("N124").Value = -45743.75
$ws.Range("H125").Value = 29966.666
$ws.Range("J125").Value = 29966.666
$ws.Range("L125").Value = 29966.666
$ws.Range("N125").Value = -39806.666
$ws.Range("H130").Value = 29600
$ws.Range("J130").Value = 29600
$ws.Range("L130").Value = 29600
$ws.Range("N130").Value = -39640
$ws.Range("H132").Value = 500030000
$ws.Range("J132").Value = 500030000
$ws.Range("L132").Value = 500030000
$ws.Range("N132").Value = -500040120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25368
$ws.Range("H92").Value = 25200.334
$ws.Range("J92").Value = 25200.334
$ws.Range("L92").Value = 25200.334
$ws.Range("N92").Value = -30192.334
$ws.Range("H96").Value = 23000
$ws.Range("J96").Value = 23000
$ws.Range("L96").Value = 23000
$ws.Range("N96").Value = -28492
$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490
$ws.Range("H106").Value = 57500
$ws.Range("J106").Value = 57500
$ws.Range("L106").Value = 57500
$ws.Range("N106").Value = -60024
$ws.Range("H141").Value = 35825.94
$ws.Range("J141").Value = 36305.516
$ws.Range("L141").Value = 36305.516
$ws.Range("N141").Value = -46665.516

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4450.7646
$ws.Range("J75").Value = 5710.4546
$ws.Range("L75").Value = 17131.3638
$ws.Range("N75").Value = -19127.3638
$ws.Range("H78").Value = 4450.7646
$ws.Range("J78").Value = 5710.4546
$ws.Range("L78").Value = 51394.0914
$ws.Range("N78").Value = -61378.0914
$ws.Range("H113").Value = 417.07144
$ws.Range("I113").Value = 386.58823
$ws.Range("K113").Value = 1159.76469
$ws.Range("M113").Value = 1010.23531
$ws.Range("H117").Value = 369.8
$ws.Range("I117").Value = 299
$ws.Range("J117").Value = 387.5
$ws.Range("K117").Value = 897
$ws.Range("L117").Value = 1162.5
$ws.Range("M117").Value = 2545
$ws.Range("N117").Value = -8046.5
$ws.Range("H129").Value = 29757.857
$ws.Range("I129").Value = 1000
$ws.Range("J129").Value = 30603.676
$ws.Range("K129").Value = 3000
$ws.Range("L129").Value = 91811.02799999999
$ws.Range("M129").Value = 2000
$ws.Range("N129").Value = -101811.028
$ws.Range("H131").Value = 51282772
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 52036924
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 156110772
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -156120852

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 51733.75
$ws.Range("I100").Value = 60501.47
$ws.Range("J100").Value = 2050
$ws.Range("K100").Value = 60501.47
$ws.Range("L100").Value = 2050
$ws.Range("M100").Value = -59960.47
$ws.Range("N100").Value = -3132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 31950
$ws.Range("J103").Value = 31950
$ws.Range("L103").Value = 31950
$ws.Range("N103").Value = -34294
$ws.Range("H132").Value = 6646.421
$ws.Range("I132").Value = 1039.5
$ws.Range("J132").Value = 16258.286
$ws.Range("K132").Value = 3118.5
$ws.Range("L132").Value = 48774.858
$ws.Range("M132").Value = -588.5
$ws.Range("N132").Value = -53834.858
